# Add the missing results for fine grained classification Config 6
# (Count Vectorizer + TFIDF + ngram(1) + POS) on the Lucene sheet.
# Rows 57-61 already contain Model (A) / Configuration (B); this fills in
# the Precision/Recall/F1/Accuracy score columns (C:F) for the five models
# (Logistic Regression, Multinomial Naive Bayes, Support Vector Machines,
# Decision Tree, Random Forest).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lucene")

# Logistic Regression
$ws.Range("C57").Value = "0.429 0.735 0.682 0.513 0.781"
$ws.Range("D57").Value = "0.355 0.647 0.613 0.455 0.788"
$ws.Range("E57").Value = "0.278 0.671 0.549 0.352 0.647 "
$ws.Range("F57").Value = "0.870 0.765 0.823 0.876 0.955"

# Multinomial Naive Bayes
$ws.Range("C58").Value = "0.383 0.750 0.687 0.456 0.754 "
$ws.Range("D58").Value = "0.383 0.592 0.573 0.455 0.726 "
$ws.Range("E58").Value = "0.239 0.773 0.564 0.299 0.614"
$ws.Range("F58").Value = "0.879 0.743 0.810 0.877 0.948 "

# Support Vector Machines
$ws.Range("C59").Value = "0.397 0.755 0.690 0.563 0.757"
$ws.Range("D59").Value = "0.359 0.662 0.609 0.473 0.809 "
$ws.Range("E59").Value = "0.251 0.703 0.560 0.401 0.614 "
$ws.Range("F59").Value = "0.873 0.778 0.823 0.879 0.955"

# Decision Tree
$ws.Range("C60").Value = "0.298 0.670 0.398 0.235 0.654 "
$ws.Range("D60").Value = "0.216 0.525 0.531 0.355 0.808"
$ws.Range("E60").Value = "0.178 0.639 0.252 0.134 0.488"
$ws.Range("F60").Value = "0.845 0.682 0.785 0.871 0.947"

# Random Forest
$ws.Range("C61").Value = "0.244 0.763 0.555 0.329 0.709"
$ws.Range("D61").Value = "0.554 0.607 0.731 0.683 0.905"
$ws.Range("E61").Value = "0.139 0.789 0.390 0.197 0.551"
$ws.Range("F61").Value = "0.897 0.756 0.833 0.896 0.957"

# Reflect the author's final scroll position / selection on the sheet.
$excel.Goto($ws.Range("A47"), $true)
$ws.Range("C63").Select()
